$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: skos:prefLable
$ws.Range("B10").Value = "Beer ontology"

# Row 11: dct:description
$ws.Range("B11").Value = "Qualitative and quantitative variables describing beer"

# Row 12: dct:creator
$ws.Range("B12").Value = "https://www.linkedin.com/in/kristina-tomicic-6bb443108/"

# Row 13: dct:rights - license changed from CC-BY-4.0 to CC0-1.0
$ws.Range("B13").Value = "https://spdx.org/licenses/CC0-1.0"

# Row 15: pav:createdOn
$ws.Range("B15").Value = "2021-08-23T18:13+02:00"

# Row 16: pav:lastUpdatedOn
$ws.Range("B16").Value = "2021-08-23T18:13+02:00"

# Row 19: new concept beer-onto:alc_percentage
$ws.Range("A19").Value = "beer-onto:alc_percentage"
$ws.Range("B19").Value = "alc_percentage"
$ws.Range("D19").Value = "Percentage of alcohol in a unit of a beer"

# Row 20: new concept beer-onto:beer_color
$ws.Range("A20").Value = "beer-onto:beer_color"
$ws.Range("B20").Value = "beer_color"
$ws.Range("D20").Value = "Color shade of a certain beer."
